$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.424036728313183
$ws.Range("C2").Value = 0.3236541775025898
$ws.Range("D2").Value = 0.04507898464314053
$ws.Range("E2").Value = 0.0651126267579123
$ws.Range("F2").Value = 1.639132484668778
$ws.Range("M2").Value = 0.4641035925004076
$ws.Range("N2").Value = 1.535336609826786

$ws.Range("B3").Value = 1.285519530909085
$ws.Range("C3").Value = 0.283749751257119
$ws.Range("D3").Value = 0.04494209123490833
$ws.Range("E3").Value = 0.06121158810412197
$ws.Range("F3").Value = 1.574188105213778
$ws.Range("M3").Value = 0.4201255521914007
$ws.Range("N3").Value = 1.541218807427441

$ws.Range("B4").Value = 1.201211014341879
$ws.Range("C4").Value = 0.2593407130870844
$ws.Range("D4").Value = 0.04488412776229467
$ws.Range("E4").Value = 0.05886348798738084
$ws.Range("F4").Value = 1.535440081974045
$ws.Range("M4").Value = 0.3933849671949829
$ws.Range("N4").Value = 1.54540087532142

$ws.Range("B5").Value = 1.167037541836066
$ws.Range("C5").Value = 0.24941592016512
$ws.Range("D5").Value = 0.04486699105699543
$ws.Range("E5").Value = 0.05791825864116262
$ws.Range("F5").Value = 1.51993037736311
$ws.Range("M5").Value = 0.382552456019404
$ws.Range("N5").Value = 1.547247549127832

$ws.Range("B6").Value = 1.161374008881012
$ws.Range("C6").Value = 0.2477692164712835
$ws.Range("D6").Value = 0.04486453490814313
$ws.Range("E6").Value = 0.05776200080712712
$ws.Range("F6").Value = 1.517371846065373
$ws.Range("M6").Value = 0.3807575801321761
$ws.Range("N6").Value = 1.547562762902601

$ws.Range("B7").Value = 1.200749402902204
$ws.Range("C7").Value = 0.2592067759423458
$ws.Range("D7").Value = 0.04488387049747544
$ws.Range("E7").Value = 0.05885069344606109
$ws.Range("F7").Value = 1.53522978147079
$ws.Range("M7").Value = 0.3932386168815967
$ws.Range("N7").Value = 1.545425204728474

$ws.Range("B8").Value = 1.376119883514377
$ws.Range("C8").Value = 0.3098752328427565
$ws.Range("D8").Value = 0.04502632242768811
$ws.Range("E8").Value = 0.06375762959570608
$ws.Range("F8").Value = 1.616503545929433
$ws.Range("M8").Value = 0.4488846923908767
$ws.Range("N8").Value = 1.537245818491783

$ws.Range("B9").Value = 1.726085186550279
$ws.Range("C9").Value = 0.4100269203997868
$ws.Range("D9").Value = 0.04551621620715451
$ws.Range("E9").Value = 0.07376463389191201
$ws.Range("F9").Value = 1.784988902980103
$ws.Range("M9").Value = 0.5601563451744482
$ws.Range("N9").Value = 1.525776409137237

$ws.Range("B10").Value = 1.987179816296816
$ws.Range("C10").Value = 0.4841779654797165
$ws.Range("D10").Value = 0.04600951389034691
$ws.Range("E10").Value = 0.08136644222482659
$ws.Range("F10").Value = 1.914560368363624
$ws.Range("M10").Value = 0.6433245139769213
$ws.Range("N10").Value = 1.520198712223333

$ws.Range("B11").Value = 2.106884454330725
$ws.Range("C11").Value = 0.5180544894164996
$ws.Range("D11").Value = 0.04626399724232755
$ws.Range("E11").Value = 0.08488230762623772
$ws.Range("F11").Value = 1.974813923900285
$ws.Range("M11").Value = 0.6814913734585133
$ws.Range("N11").Value = 1.518293653530364

$ws.Range("B12").Value = 2.15235206318988
$ws.Range("C12").Value = 0.5309049487149764
$ws.Range("D12").Value = 0.04636478001050648
$ws.Range("E12").Value = 0.08622224657718647
$ws.Range("F12").Value = 1.997822884246318
$ws.Range("M12").Value = 0.695993890904802
$ws.Range("N12").Value = 1.51766432021094

$ws.Range("B13").Value = 2.142553593153252
$ws.Range("C13").Value = 0.528136365636044
$ws.Range("D13").Value = 0.04634287680991633
$ws.Range("E13").Value = 0.08593328213702023
$ws.Range("F13").Value = 1.992858882369234
$ws.Range("M13").Value = 0.6928682838774876
$ws.Range("N13").Value = 1.5177957449811

$ws.Range("B14").Value = 2.110622310977703
$ws.Range("C14").Value = 0.5191112527258497
$ws.Range("D14").Value = 0.04627219968785568
$ws.Range("E14").Value = 0.08499237226347134
$ws.Range("F14").Value = 1.976703010314225
$ws.Range("M14").Value = 0.68268350163234
$ws.Range("N14").Value = 1.518240026001294

$ws.Range("B15").Value = 2.091081582217782
$ws.Range("C15").Value = 0.5135860355192676
$ws.Range("D15").Value = 0.04622948562909102
$ws.Range("E15").Value = 0.08441715968757535
$ws.Range("F15").Value = 1.966832220685575
$ws.Range("M15").Value = 0.6764515312012804
$ws.Range("N15").Value = 1.518524185800743

$ws.Range("B16").Value = 1.97937592942435
$ws.Range("C16").Value = 0.4819670823923161
$ws.Range("D16").Value = 0.04599349574555589
$ws.Range("E16").Value = 0.08113785584102828
$ws.Range("F16").Value = 1.910649353383661
$ws.Range("M16").Value = 0.6408370647884425
$ws.Range("N16").Value = 1.52033603317679

$ws.Range("B17").Value = 1.911089589106041
$ws.Range("C17").Value = 0.4626079122068631
$ws.Range("D17").Value = 0.04585649130098091
$ws.Range("E17").Value = 0.07914108093479655
$ws.Range("F17").Value = 1.87652137057708
$ws.Range("M17").Value = 0.6190752488235773
$ws.Range("N17").Value = 1.521610295426342

$ws.Range("B18").Value = 1.871900506117299
$ws.Range("C18").Value = 0.4514865818217686
$ws.Range("D18").Value = 0.04578051574191022
$ws.Range("E18").Value = 0.0779980106564011
$ws.Range("F18").Value = 1.857015012658849
$ws.Range("M18").Value = 0.6065896944466971
$ws.Range("N18").Value = 1.522402659140639

$ws.Range("B19").Value = 1.858646645007582
$ws.Range("C19").Value = 0.4477233780573329
$ws.Range("D19").Value = 0.04575527433487991
$ws.Range("E19").Value = 0.07761191034523307
$ws.Range("F19").Value = 1.850431533835973
$ws.Range("M19").Value = 0.60236762099629
$ws.Range("N19").Value = 1.522681114199287

$ws.Range("B20").Value = 1.918349709972347
$ws.Range("C20").Value = 0.4646673167825384
$ws.Range("D20").Value = 0.04587078255442378
$ws.Range("E20").Value = 0.07935307809231773
$ws.Range("F20").Value = 1.880141579567464
$ws.Range("M20").Value = 0.6213885835580726
$ws.Range("N20").Value = 1.521468487744983

$ws.Range("B21").Value = 2.119997522320887
$ws.Range("C21").Value = 0.5217615351793938
$ws.Range("D21").Value = 0.04629283871240375
$ws.Range("E21").Value = 0.08526850604746272
$ws.Range("F21").Value = 1.98144313385751
$ws.Range("M21").Value = 0.6856736609248912
$ws.Range("N21").Value = 1.518107021846191

$ws.Range("B22").Value = 2.252592724921442
$ws.Range("C22").Value = 0.5592056774826233
$ws.Range("D22").Value = 0.0465944562997862
$ws.Range("E22").Value = 0.0891845765622179
$ws.Range("F22").Value = 2.048771900394655
$ws.Range("M22").Value = 0.727977308238394
$ws.Range("N22").Value = 1.516447316778013

$ws.Range("B23").Value = 2.181748969656098
$ws.Range("C23").Value = 0.5392087149083977
$ws.Range("D23").Value = 0.04643108834321197
$ws.Range("E23").Value = 0.08708983740765319
$ws.Range("F23").Value = 2.012733320925321
$ws.Range("M23").Value = 0.7053720051120109
$ws.Range("N23").Value = 1.517283593028182

$ws.Range("B24").Value = 1.915067192714673
$ws.Range("C24").Value = 0.4637362335541297
$ws.Range("D24").Value = 0.04586431279468428
$ws.Range("E24").Value = 0.07925721892691939
$ws.Range("F24").Value = 1.878504527427594
$ws.Range("M24").Value = 0.6203426452986065
$ws.Range("N24").Value = 1.521532412904037

$ws.Range("B25").Value = 1.630732611316034
$ws.Range("C25").Value = 0.3828396058995054
$ws.Range("D25").Value = 0.04536063297872772
$ws.Range("E25").Value = 0.07101474698392707
$ws.Range("F25").Value = 1.738410230205375
$ws.Range("M25").Value = 0.5298133770013607
$ws.Range("N25").Value = 1.528383477871927

